# Trade #87 closed at 2026-02-17 15:53:40 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades and MarketMaking sheets
# with the results of the newly closed trade #87 (MarketMaking / DOWN).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.92   # Current Capital
$summary.Range("B4").Value = -0.09     # Total P&L $
$summary.Range("B5").Value = -0.02     # Total P&L %
$summary.Range("B6").Value = 87        # Total Trades
$summary.Range("B8").Value = 46        # Losing Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.92   # Capital
$status.Range("D4").Value = 87      # Trades
$status.Range("E4").Value = -0.09   # P&L $
$status.Range("F4").Value = -0.08   # P&L %
$status.Range("G4").Value = 33.33   # Win Rate %

# ---------------------------------------------------------------------
# New trade row data (trade #87)
# ---------------------------------------------------------------------
$rowNum = 88
$tradeRow = @{
    A = 87
    B = "2026-02-17"
    C = "15:53:33"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.16
    G = 0.11
    H = "CLOSED"
    I = -31.25
    J = -0.05
    K = 99.92
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($rowNum, 1).Value = $tradeRow.A

    # Column B holds a plain "yyyy-mm-dd" text label (not a real date), so
    # force text formatting first to stop the date-like string being
    # auto-converted into a date serial, then drop the formatting override
    # again so the cell ends up with the default (unstyled) look, matching
    # every other row in the sheet.
    $ws.Cells.Item($rowNum, 2).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 2).Value = $tradeRow.B
    $ws.Cells.Item($rowNum, 2).ClearFormats()

    $ws.Cells.Item($rowNum, 3).Value = $tradeRow.C
    $ws.Cells.Item($rowNum, 4).Value = $tradeRow.D
    $ws.Cells.Item($rowNum, 5).Value = $tradeRow.E
    $ws.Cells.Item($rowNum, 6).Value = $tradeRow.F
    $ws.Cells.Item($rowNum, 7).Value = $tradeRow.G
    $ws.Cells.Item($rowNum, 8).Value = $tradeRow.H
    $ws.Cells.Item($rowNum, 9).Value = $tradeRow.I
    $ws.Cells.Item($rowNum, 10).Value = $tradeRow.J
    $ws.Cells.Item($rowNum, 11).Value = $tradeRow.K
    $ws.Cells.Item($rowNum, 12).Value = $tradeRow.L
    $ws.Cells.Item($rowNum, 13).Value = $tradeRow.M
    $ws.Cells.Item($rowNum, 14).Value = $tradeRow.N
    $ws.Cells.Item($rowNum, 15).Value = $tradeRow.O
    $ws.Cells.Item($rowNum, 16).Value = $tradeRow.P
    $ws.Cells.Item($rowNum, 17).Value = $tradeRow.Q
}
